# Generate Report for Handoff
#
# Updates the localization-status workbook with a fresh handoff report:
#  - Priority for the "18b6e632 / 64ff716d / 67946394 / 93fa733f" file rows
#    moves from "low" to "ht" on both the zh-cn and de-de sheets.
#  - The "Latest Handoff Datetime" for those same rows is refreshed to the
#    new generation timestamp (per-locale).
#  - The Overview sheet's "Latest HO Xliff Generate Date" for those rows
#    (status "Ready for handoff") is refreshed to match the de-de timestamp.

$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhRows = 4, 5, 6, 7
foreach ($r in $zhRows) {
    $wsZh.Range("E$r").Value = "ht"
    $wsZh.Range("H$r").Value = "2016-08-20 22:38:25"
}

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$deRows = 4, 5, 6, 7
foreach ($r in $deRows) {
    $wsDe.Range("E$r").Value = "ht"
    $wsDe.Range("H$r").Value = "2016-08-20 22:38:29"
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$ovRows = 4, 5, 6, 7
foreach ($r in $ovRows) {
    $wsOverview.Range("G$r").Value = "2016-08-20 22:38:29"
}
